$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Price (D) and Volume(1h) (E) columns for each coin row.
# Price cells are forced to Text via NumberFormat "@" so that dotted/decimal
# looking strings (e.g. "0.9971", "27.361.18") are stored as literal text
# instead of being auto-coerced into numbers, then the number format is
# reset back to the default "Normal" style so no visible formatting changes.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '27.361.18'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +1.03%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.781.19'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +4.09%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9971'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.52%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '314.04'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.98%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.9973'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.45%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5201'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +10.11%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3616'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +5.85%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '42.45'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +1.33%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07357'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +1.48%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.097'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +5.90%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.9963'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.46%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '20.67'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +4.55%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.074'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +3.90%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '1.775.34'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +3.66%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '6.997'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +2.11%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '88.59'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.52%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.00001047'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.79%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06429'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +1.32%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.9980'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.30%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '16.75'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +1.60%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.849'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +4.35%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '27.415.81'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.10%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.35'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +4.47%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.068'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -2.23%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '154.91'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.55%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '20.21'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +3.97%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.356'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +13.38%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.977.58'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +3.50%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '121.62'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +2.24%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.070'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +5.86%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.09692'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +6.01%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.594'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +5.57%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.588'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.06%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.02235'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +1.57%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.06006'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +3.45%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '11.24'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +2.84%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.2035'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +2.45%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '4.848'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +2.49%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.6142'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +4.83%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.430'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +2.28%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '8.008'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +7.59%  '

$ws.Range("E43").Value = '  +3.16%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '13.26'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +6.34%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.5775'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +2.54%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.630'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +2.01%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '121.49'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +4.04%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.891'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +3.33%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.110'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +2.81%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.06716'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +1.23%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '71.02'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +2.28%  '
